$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source workbook stores these values as text (inline strings), not
# numbers, e.g. "312.94" and "-0.07%" are literal text, not a number /
# percentage. Force each target cell to the Text number format before
# writing so the COM layer doesn't silently coerce the new value into a
# numeric/percentage type.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "312.94"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.07%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "36.84"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-2.50%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.129"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.33%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07861"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.53%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.417"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-0.17%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.386"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.84%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.870"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-2.49%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.942"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "2.80%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9286"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.77%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1230"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-0.94%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1907"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.86%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.08881"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-3.10%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03272"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-2.33%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09524"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.85%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001388"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.32%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006150"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "7.16%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.373"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-3.70%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3465"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.64%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.447"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "22.13%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1297"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.90%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04328"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-1.21%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001196"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-4.30%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004350"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.75%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001324"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "8.40%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003961"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02262"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-0.97%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05116"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.20%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007459"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.24%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "1.27%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.008482"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-3.66%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.001988"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "3.44%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.007849"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-9.03%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006326"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-6.03%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.57%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002853"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-14.90%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001681"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "39.98%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002090"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.57%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0001990"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.57%"
